$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.751636028289795
$ws.Range("B1").Value = 2.069686651229858
$ws.Range("C1").Value = 2.21508264541626
$ws.Range("D1").Value = 2.929163217544556
$ws.Range("E1").Value = 1.757104396820068
